# Change the table style on the three tables (slides 14, 15, 16) from the
# custom "Table_0" style to the built-in "No Style, No Grid" table style.
$p = $ppt.ActivePresentation

$newStyleId = "{DEB89DCD-42B7-4F4A-BB58-BEA2FD30344F}"
$slideIndexes = @(14, 15, 16)

foreach ($slideIdx in $slideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}
